$wb = $excel.ActiveWorkbook

$wsDefault = $wb.Worksheets.Item("default")
$wsBrc = $wb.Worksheets.Item("brc")

# Delete the "referee_address" row (row 43) from both sheets, shifting rows up.
$wsDefault.Rows("43:43").Delete()
$wsBrc.Rows("43:43").Delete()

# Make "default" the active sheet/tab, matching the reverted state.
$wsDefault.Activate()
$wsDefault.Select()

# Update view/selection state to match the target (scrolled near the bottom, A43 selected).
$wsDefault.Application.ActiveWindow.ScrollRow = 28
$wsDefault.Range("A43").Select()

$wsBrc.Application.ActiveWindow.ScrollRow = 25
$wsBrc.Range("A43").Select()
